# New pricing data entered and documentation on Excel was updated.
#
# "Raw Data" sheet: a new tier/offer row is inserted for Product 1
# (Tier 1, Offer 4) and the price/probability figures for several
# existing rows are refreshed with the latest pricing.
#
# "pricing" sheet: the same refreshed pricing data, plus a new trailing
# row (ID 13) for Service 2 / Tier 1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# "Raw Data" sheet
# ---------------------------------------------------------------
$wsRaw = $wb.Worksheets.Item("Raw Data")

# Insert a new row above row 5 to make room for the new Tier 1 / Offer 4
# pricing entry for Product 1 - everything below shifts down by one row.
$wsRaw.Rows.Item(5).Insert()

$rawData = @(
  @(2,  "Product 1", 1, 1, 0,      0.5),
  @(3,  $null,        $null, 2, 50,   0.25),
  @(4,  $null,        $null, 3, 250,  0.15),
  @(5,  $null,        $null, 4, 1000, 0.1),
  @(6,  $null,        2, 1, 25,    0.8),
  @(7,  $null,        $null, 2, 50,   0.1),
  @(8,  $null,        $null, 3, 100,  0.1),
  @(9,  "Product 2", 1, 1, 500000, 0.05),
  @(10, "Product 3", 1, 1, 10000,  0.3),
  @(11, $null,        2, 1, 20000,  0.7),
  @(12, "Service 1", 1, 1, 30000,  0.4),
  @(13, $null,        $null, 2, 40000,  0.6),
  @(14, "Service 2", 1, 1, 500000, 1)
)

foreach ($row in $rawData) {
  $r = $row[0]
  if ($row[1] -ne $null) { $wsRaw.Cells.Item($r, 1).Value = $row[1] }
  if ($row[2] -ne $null) { $wsRaw.Cells.Item($r, 2).Value = $row[2] }
  $wsRaw.Cells.Item($r, 3).Value = $row[3]
  $wsRaw.Cells.Item($r, 4).Value = $row[4]
  $wsRaw.Cells.Item($r, 5).Value = $row[5]
}

# ---------------------------------------------------------------
# "pricing" sheet
# ---------------------------------------------------------------
$wsPricing = $wb.Worksheets.Item("pricing")

$pricingData = @(
  @(2,  1,  "Tier 1", 1, 0,      0.5,  1),
  @(3,  2,  "Tier 1", 2, 50,     0.25, 1),
  @(4,  3,  "Tier 1", 3, 250,    0.15, 1),
  @(5,  4,  "Tier 1", 4, 1000,   0.1,  1),
  @(6,  5,  "Tier 2", 1, 25,     0.8,  1),
  @(7,  6,  "Tier 2", 2, 50,     0.1,  1),
  @(8,  7,  "Tier 2", 3, 100,    0.1,  1),
  @(9,  8,  "Tier 1", 1, 500000, 0.05, 2),
  @(10, 9,  "Tier 1", 1, 10000,  0.3,  3),
  @(11, 10, "Tier 2", 2, 20000,  0.7,  3),
  @(12, 11, "Tier 1", 1, 30000,  0.4,  4),
  @(13, 12, "Tier 1", 2, 40000,  0.6,  4),
  @(14, 13, "Tier 1", 1, 500000, 1,    5)
)

foreach ($row in $pricingData) {
  $r = $row[0]
  $wsPricing.Cells.Item($r, 1).Value = $row[1]
  $wsPricing.Cells.Item($r, 2).Value = $row[2]
  $wsPricing.Cells.Item($r, 3).Value = $row[3]
  $wsPricing.Cells.Item($r, 4).Value = $row[4]
  $wsPricing.Cells.Item($r, 5).Value = $row[5]
  $wsPricing.Cells.Item($r, 6).Value = $row[6]
}

# ---------------------------------------------------------------
# Selections / active sheet to mirror the documentation updates
# ---------------------------------------------------------------
# Set the non-active sheet's selection first, then finish on "Raw Data"
# (activating it last) so it remains the active/selected tab.
$wsPricing.Range("E10").Select()

$wsRaw.Activate()
$wsRaw.Range("E2:E5").Select()
